$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 10
$ws.Range("H10").Value = 8916.666999999999
$ws.Range("I10").Value = 6800
$ws.Range("J10").Value = 19500
$ws.Range("K10").Value = 6800
$ws.Range("L10").Value = 19500
$ws.Range("M10").Value = -6507
$ws.Range("N10").Value = -20086

# Row 51
$ws.Range("H51").Value = 2323.5
$ws.Range("I51").Value = 1150
$ws.Range("J51").Value = 2910.25
$ws.Range("K51").Value = 1150
$ws.Range("L51").Value = 2910.25
$ws.Range("M51").Value = -666
$ws.Range("N51").Value = -3878.25

# Row 132
$ws.Range("H132").Value = 268548.1
$ws.Range("I132").Value = 329183.84
$ws.Range("J132").Value = 64591.453
$ws.Range("K132").Value = 987551.52
$ws.Range("L132").Value = 193774.359
$ws.Range("M132").Value = -985021.52
$ws.Range("N132").Value = -198834.359

# Row 133
$ws.Range("H133").Value = 26326.363
$ws.Range("J133").Value = 26326.363
$ws.Range("L133").Value = 26326.363
$ws.Range("N133").Value = -36446.363

# Row 138
$ws.Range("H138").Value = 3974680.8
$ws.Range("I138").Value = 1821374.5
$ws.Range("J138").Value = 4632635.5
$ws.Range("K138").Value = 5464123.5
$ws.Range("L138").Value = 13897906.5
$ws.Range("M138").Value = -5458983.5
$ws.Range("N138").Value = -13908186.5


$ws = $wb.Worksheets.Item("ARM")
# Row 3
$ws.Range("H3").Value = 1100
$ws.Range("I3").Value = 1100
$ws.Range("K3").Value = 1100
$ws.Range("M3").Value = -985

# Row 32
$ws.Range("H32").Value = 2432.321
$ws.Range("I32").Value = 2453.4932
$ws.Range("J32").Value = 2239.125
$ws.Range("K32").Value = 2453.4932
$ws.Range("L32").Value = 2239.125
$ws.Range("M32").Value = -2166.4932
$ws.Range("N32").Value = -2813.125

# Row 74
$ws.Range("H74").Value = 8871.444
$ws.Range("I74").Value = 2484.2222
$ws.Range("K74").Value = 2484.2222
$ws.Range("M74").Value = -1610.2222

# Row 77
$ws.Range("H77").Value = 8871.444
$ws.Range("I77").Value = 2484.2222
$ws.Range("K77").Value = 12421.111
$ws.Range("M77").Value = -8053.111000000001

# Row 97
$ws.Range("H97").Value = 19614264
$ws.Range("I97").Value = 23817106
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 23817106
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -23816610
$ws.Range("N97").Value = -1992

# Row 133
$ws.Range("H133").Value = 48000
$ws.Range("J133").Value = 48000
$ws.Range("L133").Value = 48000
$ws.Range("N133").Value = -53060

# Row 139
$ws.Range("H139").Value = 42918.332
$ws.Range("J139").Value = 42918.332
$ws.Range("L139").Value = 42918.332
$ws.Range("N139").Value = -53198.332


$ws = $wb.Worksheets.Item("BSM")
# Row 46
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""

# Row 94
$ws.Range("H94").Value = 1824.2222
$ws.Range("I94").Value = 1816.5
$ws.Range("J94").Value = 1839.6666
$ws.Range("K94").Value = 1816.5
$ws.Range("L94").Value = 1839.6666
$ws.Range("M94").Value = -1365.5
$ws.Range("N94").Value = -2741.6666

# Row 133
$ws.Range("H133").Value = 43226.668
$ws.Range("J133").Value = 43226.668
$ws.Range("L133").Value = 43226.668
$ws.Range("N133").Value = -53346.668

# Row 134
$ws.Range("H134").Value = 5189.25
$ws.Range("I134").Value = 5204.8
$ws.Range("K134").Value = 15614.4
$ws.Range("M134").Value = -13079.4


$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1242.7354
$ws.Range("I31").Value = 850.36
$ws.Range("J31").Value = 2332.6667
$ws.Range("K31").Value = 850.36
$ws.Range("L31").Value = 2332.6667
$ws.Range("M31").Value = -555.36
$ws.Range("N31").Value = -2922.6667

# Row 34
$ws.Range("H34").Value = 1242.7354
$ws.Range("I34").Value = 850.36
$ws.Range("J34").Value = 2332.6667
$ws.Range("K34").Value = 850.36
$ws.Range("L34").Value = 2332.6667
$ws.Range("M34").Value = -648.36
$ws.Range("N34").Value = -2736.6667

# Row 58
$ws.Range("H58").Value = 1737.75
$ws.Range("I58").Value = 982.8182
$ws.Range("K58").Value = 982.8182
$ws.Range("M58").Value = -779.8182

# Row 68
$ws.Range("H68").Value = 33266.668
$ws.Range("J68").Value = 33266.668
$ws.Range("L68").Value = 33266.668
$ws.Range("N68").Value = -34764.668

# Row 71
$ws.Range("H71").Value = 33266.668
$ws.Range("J71").Value = 33266.668
$ws.Range("L71").Value = 99800.00399999999
$ws.Range("N71").Value = -107288.004

# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = ""

# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = ""

# Row 134
$ws.Range("H134").Value = 3194.0344
$ws.Range("I134").Value = 1701.5625
$ws.Range("K134").Value = 5104.6875
$ws.Range("M134").Value = -2569.6875

# Row 136
$ws.Range("H136").Value = 1737.75
$ws.Range("I136").Value = 982.8182
$ws.Range("K136").Value = 2948.4546
$ws.Range("M136").Value = -398.4546


$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1629.2222
$ws.Range("I5").Value = 857.3333
$ws.Range("J5").Value = 2015.1666
$ws.Range("K5").Value = 2571.9999
$ws.Range("L5").Value = 6045.4998
$ws.Range("M5").Value = -2459.9999
$ws.Range("N5").Value = -6269.4998

# Row 10
$ws.Range("H10").Value = 500
$ws.Range("I10").Value = 200
$ws.Range("J10").Value = 2000
$ws.Range("K10").Value = 600
$ws.Range("L10").Value = 6000
$ws.Range("M10").Value = -461
$ws.Range("N10").Value = -6278

# Row 131
$ws.Range("H131").Value = 1669.1428
$ws.Range("J131").Value = 1904.8292
$ws.Range("L131").Value = 5714.487599999999
$ws.Range("N131").Value = -15794.4876

# Row 135
$ws.Range("H135").Value = 1629.2222
$ws.Range("I135").Value = 857.3333
$ws.Range("J135").Value = 2015.1666
$ws.Range("K135").Value = 7715.9997
$ws.Range("L135").Value = 18136.4994
$ws.Range("M135").Value = -5180.9997
$ws.Range("N135").Value = -23206.4994


$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 1163
$ws.Range("I97").Value = 1090
$ws.Range("J97").Value = 1217.75
$ws.Range("K97").Value = 1090
$ws.Range("L97").Value = 1217.75
$ws.Range("M97").Value = -594
$ws.Range("N97").Value = -2209.75

# Row 132
$ws.Range("H132").Value = 3093.327
$ws.Range("I132").Value = 2755.7026
$ws.Range("K132").Value = 8267.1078
$ws.Range("M132").Value = -5737.1078

# Row 138
$ws.Range("H138").Value = 63000
$ws.Range("J138").Value = 63000
$ws.Range("L138").Value = 63000
$ws.Range("N138").Value = -73280


$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 3314.5
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3383.6843
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3383.6843
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3607.6843

# Row 55
$ws.Range("H55").Value = 334.3846
$ws.Range("I55").Value = 227.42857
$ws.Range("K55").Value = 227.42857
$ws.Range("M55").Value = -54.42857000000001

# Row 61
$ws.Range("H61").Value = 6415.905
$ws.Range("I61").Value = 6933.375
$ws.Range("J61").Value = 4760
$ws.Range("K61").Value = 6933.375
$ws.Range("L61").Value = 4760
$ws.Range("M61").Value = -6731.375
$ws.Range("N61").Value = -5164

# Row 113
$ws.Range("H113").Value = 6415.905
$ws.Range("I113").Value = 6933.375
$ws.Range("J113").Value = 4760
$ws.Range("K113").Value = 6933.375
$ws.Range("L113").Value = 4760
$ws.Range("M113").Value = -4763.375
$ws.Range("N113").Value = -9100

# Row 126
$ws.Range("H126").Value = 3314.5
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3383.6843
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 10151.0529
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -15091.0529

# Row 132
$ws.Range("H132").Value = 5724.3335
$ws.Range("I132").Value = 5101
$ws.Range("J132").Value = 6555.4443
$ws.Range("K132").Value = 15303
$ws.Range("L132").Value = 19666.3329
$ws.Range("M132").Value = -12773
$ws.Range("N132").Value = -24726.3329


$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 888.1579
$ws.Range("I96").Value = 594.8889
$ws.Range("J96").Value = 1152.1
$ws.Range("K96").Value = 594.8889
$ws.Range("L96").Value = 1152.1
$ws.Range("M96").Value = 778.1111
$ws.Range("N96").Value = -3898.1

# Row 128
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").Value = ""

# Row 132
$ws.Range("H132").Value = 17243706
$ws.Range("I132").Value = 21741160
$ws.Range("K132").Value = 65223480
$ws.Range("M132").Value = -65220950

